# Aktualisierter Projektbericht, Rapportierung, Planung
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- First, relocate the TOTAL row (currently row 50) down to row 57 to make ---
# --- room for the newly rapportierte Arbeiten (rows 48-52) plus 3 blank rows ---
$ws.Range("A50:C50").Copy()
$ws.Range("A57:C57").PasteSpecial(-4104)
$ws.Range("A50:C50").ClearContents()
$ws.Range("C57").Formula = "=SUM(C2:C52)"

# --- Update existing entries (B41, B43:B45) with refined task descriptions ---
# Besprechung Betreuer, Vorbereitung Zwischenreview -> Handbuch..., Projektbericht, E 4.3
$ws.Range("B41").Value = "Handbuch, Aspekt 2 Adapters and Dependency Services, Projektbericht, E 4.3"
$ws.Range("B43").Value = "Handbuch, Aspekt 2 Adapters and Dependency Services, E 4.3"
$ws.Range("B44").Value = "Handbuch, Aspekt 2 Adapters and Dependency Services, E 4.3"
$ws.Range("B45").Value = "Handbuch, Aspekt 2 Adapters and Dependency Services, E 4.3"

# --- Fill in row 48, which existed already but only had a date placeholder ---
$ws.Range("A48").Value = 41465
$ws.Range("B48").Value = "Zwischenreview und Vorbereitung "
$ws.Range("C48").Value = 2

# --- New rows 49-55: copy formatting (date number format + borders) from row 47 ---
$ws.Range("A47:C47").Copy()
$ws.Range("A49:C55").PasteSpecial(-4122)

$ws.Range("A49").Value = 41466
$ws.Range("B49").Value = "Nachbesprechung, Planung angepasst "
$ws.Range("C49").Value = 1

$ws.Range("A50").Value = 41467
$ws.Range("B50").Value = "Projektbericht angegangen"
$ws.Range("C50").Value = 8

$ws.Range("A51").Value = 41468
$ws.Range("B51").Value = "Aktualisieren Wissen Handler / Command / Menü E4"
$ws.Range("C51").Value = 1

$ws.Range("A52").Value = 41468
$ws.Range("B52").Value = "Aktualisieren Wissen Action / Handler / Command / Menü E3"
$ws.Range("C52").Value = 2
$ws.Range("D52").Formula = "=SUM(C47:C52)"

# --- Rows 53-55 stay blank separators (only the date column keeps its number format) ---
$ws.Range("B53:C55").ClearContents()

# --- Sheet view: mirror final cursor position/scroll from the authored session ---
$ws.Range("D53").Select()
$excel.ActiveWindow.ScrollRow = 19
